# Generate Report for Handback
# Update generated/handoff/handback timestamps on the report sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 7a78ced6-... (row 4, column G)
$overview.Range("G4").Value = "2016-08-31 14:57:35"

# zh-cn sheet: handoff / handback datetimes for 7a78ced6-... (row 4)
$zhcn.Range("H4").Value = "2016-08-31 14:57:30"
$zhcn.Range("K4").Value = "2016-08-31 14:57:50"

# de-de sheet: handoff datetime shares the same value as the Overview sheet,
# and the handback datetime for 7a78ced6-... (row 4)
$dede.Range("H4").Value = "2016-08-31 14:57:35"
$dede.Range("K4").Value = "2016-08-31 14:57:58"
